$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted as row 65, pushing the
# existing data (previously rows 65-163) down by one row (now 66-164).
$ws.Rows.Item(65).Insert()

$ws.Cells.Item(65, 1).Value  = 4
$ws.Cells.Item(65, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(65, 3).Value  = "Los Lagos"
$ws.Cells.Item(65, 4).Value  = 44540
$ws.Cells.Item(65, 5).Value  = 10
$ws.Cells.Item(65, 6).Value  = 100112032
$ws.Cells.Item(65, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(65, 8).Value  = "Sin especificar"
$ws.Cells.Item(65, 9).Value  = "Primera"
$ws.Cells.Item(65, 10).Value = 250
$ws.Cells.Item(65, 11).Value = 11000
$ws.Cells.Item(65, 12).Value = 11000
$ws.Cells.Item(65, 13).Value = 11000
$ws.Cells.Item(65, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(65, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(65, 16).Value = 220
$ws.Cells.Item(65, 17).Value = 50
$ws.Cells.Item(65, 18).Value = "Hortaliza"
